$wb = $excel.ActiveWorkbook

# Sheets that carry the 漫展 (convention) listing rows: "展览" (index 1) and "全部类型" (index 4).
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # New row 2 takes the data that used to live in row 3 (2024-05-18 HP event).
    # The leading apostrophe forces the date-looking string to stay plain text
    # (matching the original inlineStr cells) instead of being parsed into a
    # date serial number; resetting the style afterwards drops the resulting
    # "quote prefix" style so the cell keeps the default (no explicit s=) look.
    $ws.Range("B2").Value = "'2024-05-18"
    $ws.Range("B2").Style = "Normal"
    $ws.Range("C2").Value = "丽水·第三届HP国风动漫游戏嘉年华"
    $ws.Range("D2").Value = "好溪路与望城路交汇西北侧地块 丽水市水上运动中心"
    $ws.Range("E2").Value = "2024.05.18 09:00-05.18 17:00"
    $ws.Range("F2").Value = 197
    $ws.Range("G2").Value = 68
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=82901"
    $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202403/sl5TubQI1710410535537.jpeg"

    # New row 3 takes the data that used to live in row 4 (2024-06-01 event).
    $ws.Range("B3").Value = "'2024-06-01"
    $ws.Range("B3").Style = "Normal"
    $ws.Range("C3").Value = "丽水·动漫游戏展"
    $ws.Range("D3").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E3").Value = "2024.06.01 10:00-06.01 17:00"
    $ws.Range("F3").Value = 146
    $ws.Range("G3").Value = 45
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84450"
    $ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202404/tdhb9QSW1713333412467.jpeg"

    # Old row 4 (now duplicated into row 3) is removed entirely, shrinking the sheet
    # dimension from A1:I4 down to A1:I3.
    $ws.Rows.Item(4).Delete()
}
